# Update the workbook to reflect progress on chapter 9 (5 more problems finished)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9 (chapter index 8): Progress (column E) goes from 7 to 10
$ws.Range("E9").Value = 10

# Row 10 (chapter index 9): Progress (column E) newly recorded as 10
$ws.Range("E10").Value = 10

# Recalculate so dependent formulas (e.g. G3) refresh
$excel.CalculateFull()

# Move the active selection to E10, matching the latest edit location
$ws.Range("E10").Select()

$wb.Save()
